$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "jdIYv823"
$ws.Range("B2").Value = 231011128
$ws.Range("C2").Value = "lykndww37"
$ws.Range("D2").Value = "JS!f5g9&"
$ws.Range("F2").Value = "PYaMUbVo"
$ws.Range("G2").Value = "qSxD"

# Row 3
$ws.Range("A3").Value = "SxZgA746"
$ws.Range("B3").Value = 231011127
$ws.Range("C3").Value = "fxsbida85"
$ws.Range("D3").Value = "E3!&9avU"
$ws.Range("F3").Value = "ekZKVDaf"
$ws.Range("G3").Value = "JxDn"
